$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per the authoritative diff. Numeric-looking text values are
# written with a leading apostrophe (Excel's text-entry prefix) so they stay
# plain text (matching the source inlineStr cells) instead of being parsed
# into numbers/percentages, then the quote-prefix style flag is cleared so
# the cell keeps the workbook's default (unstyled) look.

$ws.Range("D2").Value = "'331.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.35%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.34%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.657"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.27%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08341"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.71%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.790"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.52%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.986"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-4.62%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.481"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.11%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.906"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.83%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9256"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.00%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1289"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.74%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1970"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.42%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09427"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.84%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03868"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'4.06%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.1061"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.84%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001303"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.30%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006104"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.77%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.441"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.80%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'9.044"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.57%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1362"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-3.85%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-7.26%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04410"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.69%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'1.17%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004398"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.86%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-3.23%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02826"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.57%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05506"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.07%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007806"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.36%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'0.64%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.009310"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-6.03%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'2.38%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01108"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.67%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007103"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'4.63%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.00%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003400"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'13.52%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002278"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.00%"
$ws.Range("E51").Style = "Normal"
